$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: apply the "new 2020 column" numeric/text style to a cell by
# borrowing font/border from a template cell, then tweaking number format
# and alignment to match the added cellXfs entries in the target workbook.
function Set-ColStyle($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
    $ws.Range($dstAddr).NumberFormat = "0.0"
    $ws.Range($dstAddr).HorizontalAlignment = -4152
    $ws.Range($dstAddr).VerticalAlignment = -4108
}

# Row 3 header: new "2020" column header, same style as the "2018" header (D3)
$ws.Range("D3").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Value = 2020

# Row 4 (totals row, bordered like D4)
Set-ColStyle "D4" "E4"
$ws.Range("E4").Value = 22.1

# Rows 5-18 (regular data rows, same look as D5..D23)
Set-ColStyle "D5" "E5"

Set-ColStyle "D5" "E6"
$ws.Range("E6").Value = 52.7

Set-ColStyle "D5" "E7"
$ws.Range("E7").Value = 5

Set-ColStyle "D5" "E8"

Set-ColStyle "D5" "E9"
$ws.Range("E9").Value = 4.8

Set-ColStyle "D5" "E10"
$ws.Range("E10").Value = 15.8

Set-ColStyle "D5" "E11"
$ws.Range("E11").Value = 13.5

Set-ColStyle "D5" "E12"
$ws.Range("E12").Value = 9.6

Set-ColStyle "D5" "E13"
$ws.Range("E13").Value = 2.7

Set-ColStyle "D5" "E14"
$ws.Range("E14").Value = 14.7

Set-ColStyle "D5" "E15"
$ws.Range("E15").Value = 18.2

Set-ColStyle "D5" "E16"
$ws.Range("E16").Value = 74

Set-ColStyle "D5" "E17"
$ws.Range("E17").Value = 35.1

Set-ColStyle "D5" "E18"

# Rows 19-23: no data available, "-" placeholder text
Set-ColStyle "D5" "E19"
$ws.Range("E19").Value = "-"

Set-ColStyle "D5" "E20"
$ws.Range("E20").Value = "-"

Set-ColStyle "D5" "E21"
$ws.Range("E21").Value = "-"

Set-ColStyle "D5" "E22"
$ws.Range("E22").Value = "-"

Set-ColStyle "D5" "E23"
$ws.Range("E23").Value = "-"

# Row 24 (sub-header row), stays empty
Set-ColStyle "D5" "E24"

# Rows 25-28: "-" placeholder text
Set-ColStyle "D5" "E25"
$ws.Range("E25").Value = "-"

Set-ColStyle "D5" "E26"
$ws.Range("E26").Value = "-"

Set-ColStyle "D5" "E27"
$ws.Range("E27").Value = "-"

Set-ColStyle "D5" "E28"
$ws.Range("E28").Value = "-"

# Row 29 (bottom, thick-bottom bordered like D29), "-" placeholder text
Set-ColStyle "D29" "E29"
$ws.Range("E29").Value = "-"

# Selection left on J24, matching the saved sheet view in the target file
$ws.Range("J24").Select() | Out-Null
